$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10
$ws.Range("E10").Value = 8.380000000000001
$ws.Range("G10").Value = 468.62
$ws.Range("H10").Value = 9.038038375901294
$ws.Range("I10").Value = 75.1961981462664
$ws.Range("K10").Value = 0.19

# Row 11
$ws.Range("E11").Value = 8.44
$ws.Range("G11").Value = 475.46
$ws.Range("H11").Value = 9.963588420820088
$ws.Range("I11").Value = 75.73501523426967
$ws.Range("K11").Value = 0.2

# Row 12
$ws.Range("E12").Value = 8.51
$ws.Range("G12").Value = 485.92
$ws.Range("H12").Value = 10.55336407969622
$ws.Range("I12").Value = 76.07276696851476
$ws.Range("K12").Value = 0.2

# Row 13
$ws.Range("E13").Value = 7.32
$ws.Range("G13").Value = 412.26
$ws.Range("H13").Value = 9.793185203867267
$ws.Range("I13").Value = 75.70516509923618
$ws.Range("K13").Value = 0.17

# Row 14
$ws.Range("E14").Value = 4.62
$ws.Range("G14").Value = 251.4
$ws.Range("H14").Value = 7.304037140469404
$ws.Range("I14").Value = 74.42060444317758
$ws.Range("K14").Value = 0.11

# Row 15
$ws.Range("G15").Value = 44.71
$ws.Range("H15").Value = 4.018826999725093
$ws.Range("I15").Value = 72.74050571763655
$ws.Range("K15").Value = 0.02

# Row 32
$ws.Range("E32").Value = 38.09999999999999
$ws.Range("G32").Value = 2138.37
$ws.Range("H32").Value = 69.97104022047938
$ws.Range("I32").Value = 449.8702556091012
$ws.Range("K32").Value = 0.8900000000000001
